$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5033.2666
$ws.Range("J40").Value = 5033.2666
$ws.Range("L40").Value = 5033.2666
$ws.Range("N40").Value = -5383.2666
$ws.Range("H42").Value = 205.11111
$ws.Range("I42").Value = 76.5
$ws.Range("K42").Value = 229.5
$ws.Range("M42").Value = 0.5
$ws.Range("H51").Value = 4999.6313
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 4999.6313
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 4999.6313
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -5967.6313
$ws.Range("H62").Value = 5187.4546
$ws.Range("I62").Value = 1350
$ws.Range("K62").Value = 1350
$ws.Range("M62").Value = -726
$ws.Range("H65").Value = 5187.4546
$ws.Range("I65").Value = 1350
$ws.Range("K65").Value = 6750
$ws.Range("M65").Value = -3630
$ws.Range("H76").Value = 15544.728
$ws.Range("I76").Value = 52501
$ws.Range("K76").Value = 52501
$ws.Range("M76").Value = -52186
$ws.Range("H79").Value = 15544.728
$ws.Range("I79").Value = 52501
$ws.Range("K79").Value = 52501
$ws.Range("M79").Value = -51409
$ws.Range("H86").Value = 2991.484
$ws.Range("I86").Value = 2344.25
$ws.Range("K86").Value = 2344.25
$ws.Range("M86").Value = -1221.25
$ws.Range("H89").Value = 2991.484
$ws.Range("I89").Value = 2344.25
$ws.Range("K89").Value = 11721.25
$ws.Range("M89").Value = -6105.25
$ws.Range("H99").Value = 589.3077
$ws.Range("I99").Value = 538.4167
$ws.Range("J99").Value = 1200
$ws.Range("K99").Value = 1615.2501
$ws.Range("L99").Value = 3600
$ws.Range("M99").Value = -117.2501
$ws.Range("N99").Value = -6596
$ws.Range("H113").Value = 4829.952
$ws.Range("I113").Value = 2665.6667
$ws.Range("K113").Value = 2665.6667
$ws.Range("M113").Value = 588.3332999999998
$ws.Range("H129").Value = 27530.625
$ws.Range("J129").Value = 31320.715
$ws.Range("L129").Value = 93962.145
$ws.Range("N129").Value = -103962.145
$ws.Range("H137").Value = 2267.3225
$ws.Range("I137").Value = 2410.2632
$ws.Range("K137").Value = 7230.7896
$ws.Range("M137").Value = -4680.7896

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1858.5714
$ws.Range("I45").Value = 1596.3077
$ws.Range("J45").Value = 2284.75
$ws.Range("K45").Value = 1596.3077
$ws.Range("L45").Value = 2284.75
$ws.Range("M45").Value = -1219.3077
$ws.Range("N45").Value = -3038.75
$ws.Range("H61").Value = 20005380
$ws.Range("I61").Value = 27781922
$ws.Range("K61").Value = 27781922
$ws.Range("M61").Value = -27781710
$ws.Range("H122").Value = 2722.5405
$ws.Range("I122").Value = 1829.3914
$ws.Range("K122").Value = 5488.174199999999
$ws.Range("M122").Value = -3038.174199999999
$ws.Range("H132").Value = 62603988
$ws.Range("I132").Value = 24751.1
$ws.Range("K132").Value = 74253.29999999999
$ws.Range("M132").Value = -71723.29999999999
$ws.Range("H135").Value = 39199.4
$ws.Range("J135").Value = 39199.4
$ws.Range("L135").Value = 39199.4
$ws.Range("N135").Value = -49339.4
$ws.Range("H136").Value = 20005380
$ws.Range("I136").Value = 27781922
$ws.Range("K136").Value = 83345766
$ws.Range("M136").Value = -83343216

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5181.364
$ws.Range("I99").Value = 3997.8572
$ws.Range("K99").Value = 3997.8572
$ws.Range("M99").Value = -2499.8572

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 38367.2
$ws.Range("J9").Value = 38367.2
$ws.Range("L9").Value = 38367.2
$ws.Range("N9").Value = -38703.2
$ws.Range("H31").Value = 4102.171
$ws.Range("I31").Value = 2845.6538
$ws.Range("J31").Value = 6280.1333
$ws.Range("K31").Value = 2845.6538
$ws.Range("L31").Value = 6280.1333
$ws.Range("M31").Value = -2550.6538
$ws.Range("N31").Value = -6870.1333
$ws.Range("H34").Value = 4102.171
$ws.Range("I34").Value = 2845.6538
$ws.Range("J34").Value = 6280.1333
$ws.Range("K34").Value = 2845.6538
$ws.Range("L34").Value = 6280.1333
$ws.Range("M34").Value = -2643.6538
$ws.Range("N34").Value = -6684.1333
$ws.Range("H99").Value = 11628.286
$ws.Range("I99").Value = 13066.333
$ws.Range("K99").Value = 13066.333
$ws.Range("M99").Value = -11568.333
$ws.Range("H105").Value = 22804.4
$ws.Range("I105").Value = 2375
$ws.Range("K105").Value = 2375
$ws.Range("M105").Value = -628
$ws.Range("H126").Value = 11628.286
$ws.Range("I126").Value = 13066.333
$ws.Range("K126").Value = 39198.999
$ws.Range("M126").Value = -36728.999
$ws.Range("H132").Value = 66114.84
$ws.Range("I132").Value = 74345.57000000001
$ws.Range("J132").Value = 8499.75
$ws.Range("K132").Value = 223036.71
$ws.Range("L132").Value = 25499.25
$ws.Range("M132").Value = -220506.71
$ws.Range("N132").Value = -30559.25

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 174347.62
$ws.Range("J37").Value = 174347.62
$ws.Range("L37").Value = 523042.86
$ws.Range("N37").Value = -523266.86
$ws.Range("H88").Value = 9313
$ws.Range("J88").Value = 9313
$ws.Range("L88").Value = 27939
$ws.Range("N88").Value = -28795
$ws.Range("H91").Value = 9313
$ws.Range("J91").Value = 9313
$ws.Range("L91").Value = 27939
$ws.Range("N91").Value = -30903
$ws.Range("H129").Value = 1921.4828
$ws.Range("I129").Value = 905.8
$ws.Range("J129").Value = 2133.0833
$ws.Range("K129").Value = 2717.4
$ws.Range("L129").Value = 6399.249899999999
$ws.Range("M129").Value = 2282.6
$ws.Range("N129").Value = -16399.2499
$ws.Range("H131").Value = 1619.0238
$ws.Range("I131").Value = 1208.1666
$ws.Range("J131").Value = 1687.5
$ws.Range("K131").Value = 3624.4998
$ws.Range("L131").Value = 5062.5
$ws.Range("M131").Value = 1415.5002
$ws.Range("N131").Value = -15142.5
$ws.Range("H133").Value = 3289.4443
$ws.Range("I133").Value = 3289.4443
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 9868.332900000001
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -4808.332900000001
$ws.Range("N133").ClearContents()

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H58").Value = 40000
$ws.Range("I58").Value = 40000
$ws.Range("K58").Value = 40000
$ws.Range("M58").Value = -39723
$ws.Range("H97").Value = 3817.353
$ws.Range("I97").Value = 3498.9
$ws.Range("J97").Value = 4272.2856
$ws.Range("K97").Value = 3498.9
$ws.Range("L97").Value = 4272.2856
$ws.Range("M97").Value = -3002.9
$ws.Range("N97").Value = -5264.2856
$ws.Range("H102").Value = 2976.6924
$ws.Range("I102").Value = 1776.909
$ws.Range("J102").Value = 3856.5334
$ws.Range("K102").Value = 1776.909
$ws.Range("L102").Value = 3856.5334
$ws.Range("M102").Value = -154.9090000000001
$ws.Range("N102").Value = -7100.5334
$ws.Range("H123").Value = 27427.572
$ws.Range("J123").Value = 9997.666999999999
$ws.Range("L123").Value = 9997.666999999999
$ws.Range("N123").Value = -14897.667
$ws.Range("H126").Value = 4714.143
$ws.Range("I126").Value = 3999.5
$ws.Range("K126").Value = 11998.5
$ws.Range("M126").Value = -9528.5
$ws.Range("H132").Value = 1305
$ws.Range("I132").Value = 1600
$ws.Range("J132").Value = 1010
$ws.Range("K132").Value = 4800
$ws.Range("L132").Value = 3030
$ws.Range("M132").Value = -2270
$ws.Range("N132").Value = -8090

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3087.9722
$ws.Range("I7").Value = 2231.611
$ws.Range("K7").Value = 2231.611
$ws.Range("M7").Value = -2119.611
$ws.Range("H40").Value = 3909.0908
$ws.Range("I40").Value = 3800
$ws.Range("K40").Value = 3800
$ws.Range("M40").Value = -3664
$ws.Range("H126").Value = 3087.9722
$ws.Range("I126").Value = 2231.611
$ws.Range("K126").Value = 6694.833
$ws.Range("M126").Value = -4224.833
